$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 213.16667
$ws.Range("I12").Value = 194.5
$ws.Range("J12").Value = 250.5
$ws.Range("K12").Value = 194.5
$ws.Range("L12").Value = 250.5
$ws.Range("M12").Value = -24.5
$ws.Range("N12").Value = -590.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 426.23077
$ws.Range("I18").Value = 367.36365
$ws.Range("J18").Value = 750
$ws.Range("K18").Value = 367.36365
$ws.Range("L18").Value = 750
$ws.Range("M18").Value = -83.36365000000001
$ws.Range("N18").Value = -1318

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 688.1429000000001
$ws.Range("I33").Value = 752.6667
$ws.Range("J33").Value = 301
$ws.Range("K33").Value = 752.6667
$ws.Range("L33").Value = 301
$ws.Range("M33").Value = -523.6667
$ws.Range("N33").Value = -759

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 10423.125
$ws.Range("I62").Value = 2665
$ws.Range("J62").Value = 15078
$ws.Range("K62").Value = 2665
$ws.Range("L62").Value = 15078
$ws.Range("M62").Value = -2041
$ws.Range("N62").Value = -16326

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 10423.125
$ws.Range("I65").Value = 2665
$ws.Range("J65").Value = 15078
$ws.Range("K65").Value = 13325
$ws.Range("L65").Value = 75390
$ws.Range("M65").Value = -10205
$ws.Range("N65").Value = -81630

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 970.8889
$ws.Range("I96").Value = 576
$ws.Range("J96").Value = 1464.5
$ws.Range("K96").Value = 1728
$ws.Range("L96").Value = 4393.5
$ws.Range("M96").Value = -355
$ws.Range("N96").Value = -7139.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H101").Value = 1179.6154
$ws.Range("I101").Value = 824.375
$ws.Range("J101").Value = 1748
$ws.Range("K101").Value = 2473.125
$ws.Range("L101").Value = 5244
$ws.Range("M101").Value = -851.125
$ws.Range("N101").Value = -8488

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1603763.8
$ws.Range("I137").Value = 2689004.2
$ws.Range("J137").Value = 1742.1428
$ws.Range("K137").Value = 8067012.600000001
$ws.Range("L137").Value = 5226.428400000001
$ws.Range("M137").Value = -8064462.600000001
$ws.Range("N137").Value = -10326.4284

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4176.527
$ws.Range("I138").Value = 6118.3335
$ws.Range("J138").Value = 3907.6616
$ws.Range("K138").Value = 18355.0005
$ws.Range("L138").Value = 11722.9848
$ws.Range("M138").Value = -13215.0005
$ws.Range("N138").Value = -22002.9848

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 197.5
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 197.5
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 197.5
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -429.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 233.33333
$ws.Range("I5").Value = 200
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 200
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = -88
$ws.Range("N5").Value = -524

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6677178.5
$ws.Range("I32").Value = 7471453.5
$ws.Range("J32").Value = 25125
$ws.Range("K32").Value = 7471453.5
$ws.Range("L32").Value = 25125
$ws.Range("M32").Value = -7471166.5
$ws.Range("N32").Value = -25699

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8549844
$ws.Range("I61").Value = 12347397
$ws.Range("J61").Value = 5351
$ws.Range("K61").Value = 12347397
$ws.Range("L61").Value = 5351
$ws.Range("M61").Value = -12347185
$ws.Range("N61").Value = -5775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 54290.156
$ws.Range("I122").Value = 68033.266
$ws.Range("J122").Value = 2753.5
$ws.Range("K122").Value = 204099.798
$ws.Range("L122").Value = 8260.5
$ws.Range("M122").Value = -201649.798
$ws.Range("N122").Value = -13160.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8549844
$ws.Range("I136").Value = 12347397
$ws.Range("J136").Value = 5351
$ws.Range("K136").Value = 37042191
$ws.Range("L136").Value = 16053
$ws.Range("M136").Value = -37039641
$ws.Range("N136").Value = -21153

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 233.33333
$ws.Range("I4").Value = 200
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 200
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -85
$ws.Range("N4").Value = -530

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 933.3333
$ws.Range("I22").Value = 933.3333
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 933.3333
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -760.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1686.2142
$ws.Range("I80").Value = 3805.2
$ws.Range("J80").Value = 509
$ws.Range("K80").Value = 3805.2
$ws.Range("L80").Value = 509
$ws.Range("M80").Value = -2807.2
$ws.Range("N80").Value = -2505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 1686.2142
$ws.Range("I83").Value = 3805.2
$ws.Range("J83").Value = 509
$ws.Range("K83").Value = 19026
$ws.Range("L83").Value = 2545
$ws.Range("M83").Value = -14034
$ws.Range("N83").Value = -12529

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2504.457
$ws.Range("I134").Value = 2413.818
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 7241.454000000001
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4706.454000000001
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.1
$ws.Range("I7").Value = 43.25
$ws.Range("J7").Value = 71.333336
$ws.Range("K7").Value = 43.25
$ws.Range("L7").Value = 71.333336
$ws.Range("M7").Value = 69.75
$ws.Range("N7").Value = -297.333336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 344.4091
$ws.Range("I22").Value = 268.23077
$ws.Range("J22").Value = 454.44446
$ws.Range("K22").Value = 268.23077
$ws.Range("L22").Value = 454.44446
$ws.Range("M22").Value = 81.76922999999999
$ws.Range("N22").Value = -1154.44446

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1626.3572
$ws.Range("I58").Value = 1483.3334
$ws.Range("J58").Value = 1733.625
$ws.Range("K58").Value = 1483.3334
$ws.Range("L58").Value = 1733.625
$ws.Range("M58").Value = -1280.3334
$ws.Range("N58").Value = -2139.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2185.6177
$ws.Range("I99").Value = 1364
$ws.Range("J99").Value = 2438.423
$ws.Range("K99").Value = 1364
$ws.Range("L99").Value = 2438.423
$ws.Range("M99").Value = 134
$ws.Range("N99").Value = -5434.423

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 2185.6177
$ws.Range("I126").Value = 1364
$ws.Range("J126").Value = 2438.423
$ws.Range("K126").Value = 4092
$ws.Range("L126").Value = 7315.268999999999
$ws.Range("M126").Value = -1622
$ws.Range("N126").Value = -12255.269

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H129").Value = 56780
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 56780
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 56780
$ws.Range("N129").Value = -66780

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1626.3572
$ws.Range("I136").Value = 1483.3334
$ws.Range("J136").Value = 1733.625
$ws.Range("K136").Value = 4450.0002
$ws.Range("L136").Value = 5200.875
$ws.Range("M136").Value = -1900.0002
$ws.Range("N136").Value = -10300.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 88235450
$ws.Range("I23").Value = 150
$ws.Range("J23").Value = 125000160
$ws.Range("K23").Value = 450
$ws.Range("L23").Value = 375000480
$ws.Range("M23").Value = -215
$ws.Range("N23").Value = -375000950

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 978975.5600000001
$ws.Range("I129").Value = 581
$ws.Range("J129").Value = 1167128.4
$ws.Range("K129").Value = 1743
$ws.Range("L129").Value = 3501385.2
$ws.Range("M129").Value = 3257
$ws.Range("N129").Value = -3511385.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4537.6665
$ws.Range("I131").Value = 525
$ws.Range("J131").Value = 5429.3706
$ws.Range("K131").Value = 1575
$ws.Range("L131").Value = 16288.1118
$ws.Range("M131").Value = 3465
$ws.Range("N131").Value = -26368.1118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 158.125
$ws.Range("I2").Value = 156.23077
$ws.Range("J2").Value = 166.33333
$ws.Range("K2").Value = 156.23077
$ws.Range("L2").Value = 166.33333
$ws.Range("M2").Value = -43.23077000000001
$ws.Range("N2").Value = -392.33333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5299.5864
$ws.Range("I70").Value = 5147.048
$ws.Range("J70").Value = 5700
$ws.Range("K70").Value = 5147.048
$ws.Range("L70").Value = 5700
$ws.Range("M70").Value = -4877.048
$ws.Range("N70").Value = -6240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5299.5864
$ws.Range("I73").Value = 5147.048
$ws.Range("J73").Value = 5700
$ws.Range("K73").Value = 5147.048
$ws.Range("L73").Value = 5700
$ws.Range("M73").Value = -4211.048
$ws.Range("N73").Value = -7572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2900
$ws.Range("I122").Value = 2900
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 8700
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -6250
$ws.Range("N122").Value = -13600

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 8479.261
$ws.Range("I123").Value = 3000
$ws.Range("J123").Value = 9632.789000000001
$ws.Range("K123").Value = 3000
$ws.Range("L123").Value = 9632.789000000001
$ws.Range("M123").Value = -550
$ws.Range("N123").Value = -14532.789

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 28576342
$ws.Range("I132").Value = 35719320
$ws.Range("J132").Value = 4430
$ws.Range("K132").Value = 107157960
$ws.Range("L132").Value = 13290
$ws.Range("M132").Value = -107155430
$ws.Range("N132").Value = -18350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 732.1875
$ws.Range("I16").Value = 265.42856
$ws.Range("J16").Value = 3999.5
$ws.Range("K16").Value = 265.42856
$ws.Range("L16").Value = 3999.5
$ws.Range("M16").Value = -95.42856
$ws.Range("N16").Value = -4339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 17223.5
$ws.Range("I22").Value = 667.6
$ws.Range("J22").Value = 100003
$ws.Range("K22").Value = 667.6
$ws.Range("L22").Value = 100003
$ws.Range("M22").Value = -372.6
$ws.Range("N22").Value = -100593

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 17223.5
$ws.Range("I27").Value = 667.6
$ws.Range("J27").Value = 100003
$ws.Range("K27").Value = 667.6
$ws.Range("L27").Value = 100003
$ws.Range("M27").Value = -560.6
$ws.Range("N27").Value = -100217

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 347.76923
$ws.Range("I55").Value = 188.33333
$ws.Range("J55").Value = 484.42856
$ws.Range("K55").Value = 188.33333
$ws.Range("L55").Value = 484.42856
$ws.Range("M55").Value = -15.33332999999999
$ws.Range("N55").Value = -830.4285600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5803.0513
$ws.Range("I122").Value = 4556
$ws.Range("J122").Value = 5986.4414
$ws.Range("K122").Value = 13668
$ws.Range("L122").Value = 17959.3242
$ws.Range("M122").Value = -11218
$ws.Range("N122").Value = -22859.3242

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3314.739
$ws.Range("I132").Value = 2384.6924
$ws.Range("J132").Value = 4523.8
$ws.Range("K132").Value = 7154.0772
$ws.Range("L132").Value = 13571.4
$ws.Range("M132").Value = -4624.0772
$ws.Range("N132").Value = -18631.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2618.3225
$ws.Range("I122").Value = 1803.6316
$ws.Range("J122").Value = 3908.25
$ws.Range("K122").Value = 5410.8948
$ws.Range("L122").Value = 11724.75
$ws.Range("M122").Value = -2960.8948
$ws.Range("N122").Value = -16624.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7885172
$ws.Range("I132").Value = 2161.7932
$ws.Range("J132").Value = 36461084
$ws.Range("K132").Value = 6485.3796
$ws.Range("L132").Value = 109383252
$ws.Range("M132").Value = -3955.3796
$ws.Range("N132").Value = -109388312
